# Update the "Priority" (column E) values for several groups of rows in
# the "data" worksheet. Each group corresponds to a contiguous block of
# rows sharing the same "From" (column B) value; the whole block's
# Priority is being reset to a new constant value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$groups = @(
    @{ Start = 87;  End = 103; Value = 0 },
    @{ Start = 138; End = 154; Value = 0 },
    @{ Start = 155; End = 171; Value = 0 },
    @{ Start = 223; End = 239; Value = 0 },
    @{ Start = 240; End = 256; Value = 2 },
    @{ Start = 257; End = 273; Value = 2 },
    @{ Start = 274; End = 290; Value = 2 }
)

foreach ($group in $groups) {
    $range = $ws.Range("E$($group.Start):E$($group.End)")
    $range.Value = $group.Value
}
